$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L225").Value = 16.07
$ws.Range("I226").Value = 445.94
$ws.Range("L226").Value = 409.19
$ws.Range("I228").Value = 744.02
$ws.Range("I229").Value = 792.87
$ws.Range("L229").Value = 768.64
$ws.Range("H230").Value = 513.35
$ws.Range("I230").Value = 810.58
$ws.Range("K230").Value = 513.35
$ws.Range("L230").Value = 781.95
$ws.Range("M230").Value = 106.44
$ws.Range("H231").Value = 495.15
$ws.Range("I231").Value = 801.96
$ws.Range("H232").Value = 424.38
$ws.Range("I232").Value = 764.6799999999999
$ws.Range("K232").Value = 424.38
$ws.Range("L232").Value = 745.63
$ws.Range("I233").Value = 687.15
$ws.Range("J233").Value = 79.48
$ws.Range("L233").Value = 669.01
$ws.Range("H234").Value = 165.07
$ws.Range("I234").Value = 533.08
$ws.Range("K234").Value = 165.07
$ws.Range("L234").Value = 509.24
$ws.Range("M234").Value = 55.75
$ws.Range("I235").Value = 177.07
$ws.Range("L235").Value = 115
$ws.Range("I249").Value = 69.28
$ws.Range("L249").Value = 26.84
$ws.Range("H250").Value = 119.58
$ws.Range("I250").Value = 454.9
$ws.Range("K250").Value = 119.58
$ws.Range("L250").Value = 420.16
$ws.Range("H251").Value = 268.31
$ws.Range("I251").Value = 651.85
$ws.Range("J251").Value = 74.72
$ws.Range("K251").Value = 268.31
$ws.Range("L251").Value = 632.89
$ws.Range("H252").Value = 395.86
$ws.Range("I252").Value = 747.33
$ws.Range("K252").Value = 395.86
$ws.Range("L252").Value = 730.5
$ws.Range("I253").Value = 795.59
$ws.Range("J253").Value = 96.06999999999999
$ws.Range("J254").Value = 98.98999999999999
$ws.Range("H255").Value = 500.98
$ws.Range("J255").Value = 97.53
$ws.Range("K255").Value = 500.98
$ws.Range("M256").Value = 93.23
$ws.Range("I257").Value = 690.98
$ws.Range("I258").Value = 539.42
$ws.Range("K258").Value = 169.49
$ws.Range("L258").Value = 515.48
$ws.Range("M258").Value = 56.54
$ws.Range("I259").Value = 187.08
$ws.Range("K259").Value = 32.6
$ws.Range("L259").Value = 120.12
$ws.Range("M259").Value = 20.05
